$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# ALC row 5 (key=5503)
$ws_ALC.Range("H5").Value = 701.4666999999999
$ws_ALC.Range("I5").Value = 724.7692
$ws_ALC.Range("K5").Value = 724.7692
$ws_ALC.Range("M5").Value = -609.7692

# ALC row 12 (key=5515)
$ws_ALC.Range("H12").Value = 999
$ws_ALC.Range("I12").Value = 0
$ws_ALC.Range("J12").Value = 999
$ws_ALC.Range("K12").Value = 0
$ws_ALC.Range("L12").Value = 999
$ws_ALC.Range("M12").Value = ""
$ws_ALC.Range("N12").Value = -1339

# ALC row 96 (key=19894)
$ws_ALC.Range("H96").Value = 1117.8334
$ws_ALC.Range("I96").Value = 696
$ws_ALC.Range("K96").Value = 2088
$ws_ALC.Range("M96").Value = -715

# ALC row 100 (key=19906)
$ws_ALC.Range("H100").Value = 1814.8889
$ws_ALC.Range("I100").Value = 1868.7059
$ws_ALC.Range("J100").Value = 900
$ws_ALC.Range("K100").Value = 1868.7059
$ws_ALC.Range("L100").Value = 900
$ws_ALC.Range("M100").Value = -1327.7059
$ws_ALC.Range("N100").Value = -1982

# ALC row 101 (key=19884)
$ws_ALC.Range("H101").Value = 673.4545000000001
$ws_ALC.Range("I101").Value = 673.4545000000001
$ws_ALC.Range("J101").Value = 0
$ws_ALC.Range("K101").Value = 2020.3635
$ws_ALC.Range("L101").Value = 0
$ws_ALC.Range("M101").Value = -398.3635000000002
$ws_ALC.Range("N101").Value = ""

# ALC row 106 (key=19903)
$ws_ALC.Range("H106").Value = 3407
$ws_ALC.Range("I106").Value = 3476.3076
$ws_ALC.Range("K106").Value = 3476.3076
$ws_ALC.Range("M106").Value = -2845.3076

# ALC row 116 (key=27778)
$ws_ALC.Range("H116").Value = 6867.6
$ws_ALC.Range("J116").Value = 8097.2856
$ws_ALC.Range("L116").Value = 8097.2856
$ws_ALC.Range("N116").Value = -14981.2856

# ALC row 137 (key=44013)
$ws_ALC.Range("H137").Value = 3052
$ws_ALC.Range("I137").Value = 1809.3077
$ws_ALC.Range("J137").Value = 3949.5
$ws_ALC.Range("K137").Value = 5427.9231
$ws_ALC.Range("L137").Value = 11848.5
$ws_ALC.Range("M137").Value = -2877.9231
$ws_ALC.Range("N137").Value = -16948.5

# ALC row 138 (key=44169)
$ws_ALC.Range("H138").Value = 5440.9
$ws_ALC.Range("I138").Value = 3931.6
$ws_ALC.Range("K138").Value = 11794.8
$ws_ALC.Range("M138").Value = -6654.799999999999

# ALC row 140 (key=42459)
$ws_ALC.Range("H140").Value = 92424.2
$ws_ALC.Range("J140").Value = 92424.2
$ws_ALC.Range("L140").Value = 92424.2
$ws_ALC.Range("N140").Value = -102784.2

$ws_ARM = $wb.Worksheets.Item("ARM")
# ARM row 61 (key=43999)
$ws_ARM.Range("H61").Value = 4883.9
$ws_ARM.Range("I61").Value = 1793.5
$ws_ARM.Range("K61").Value = 1793.5
$ws_ARM.Range("M61").Value = -1581.5

# ARM row 88 (key=12530)
$ws_ARM.Range("H88").Value = 2762.375
$ws_ARM.Range("I88").Value = 1000
$ws_ARM.Range("J88").Value = 3014.1428
$ws_ARM.Range("K88").Value = 1000
$ws_ARM.Range("L88").Value = 3014.1428
$ws_ARM.Range("M88").Value = -594
$ws_ARM.Range("N88").Value = -3826.1428

# ARM row 91 (key=12530)
$ws_ARM.Range("H91").Value = 2762.375
$ws_ARM.Range("I91").Value = 1000
$ws_ARM.Range("J91").Value = 3014.1428
$ws_ARM.Range("K91").Value = 1000
$ws_ARM.Range("L91").Value = 3014.1428
$ws_ARM.Range("M91").Value = 404
$ws_ARM.Range("N91").Value = -5822.1428

# ARM row 97 (key=19941)
$ws_ARM.Range("H97").Value = 457.0625
$ws_ARM.Range("I97").Value = 365.30768
$ws_ARM.Range("J97").Value = 854.6667
$ws_ARM.Range("K97").Value = 365.30768
$ws_ARM.Range("L97").Value = 854.6667
$ws_ARM.Range("M97").Value = 130.69232
$ws_ARM.Range("N97").Value = -1846.6667

# ARM row 136 (key=43999)
$ws_ARM.Range("H136").Value = 4883.9
$ws_ARM.Range("I136").Value = 1793.5
$ws_ARM.Range("K136").Value = 5380.5
$ws_ARM.Range("M136").Value = -2830.5

# ARM row 139 (key=42321)
$ws_ARM.Range("H139").Value = 95213.75
$ws_ARM.Range("I139").Value = 78999
$ws_ARM.Range("J139").Value = 97530.14
$ws_ARM.Range("K139").Value = 78999
$ws_ARM.Range("L139").Value = 97530.14
$ws_ARM.Range("M139").Value = -73859
$ws_ARM.Range("N139").Value = -107810.14

$ws_BSM = $wb.Worksheets.Item("BSM")
# BSM row 99 (key=19943)
$ws_BSM.Range("H99").Value = 2333.2222
$ws_BSM.Range("I99").Value = 1999.875
$ws_BSM.Range("K99").Value = 1999.875
$ws_BSM.Range("M99").Value = -501.875

# BSM row 134 (key=43998)
$ws_BSM.Range("H134").Value = 2926.9412
$ws_BSM.Range("I134").Value = 1984.875
$ws_BSM.Range("K134").Value = 5954.625
$ws_BSM.Range("M134").Value = -3419.625

$ws_CRP = $wb.Worksheets.Item("CRP")
# CRP row 31 (key=44023)
$ws_CRP.Range("H31").Value = 7161.357
$ws_CRP.Range("I31").Value = 3592.5557
$ws_CRP.Range("K31").Value = 3592.5557
$ws_CRP.Range("M31").Value = -3297.5557

# CRP row 34 (key=44023)
$ws_CRP.Range("H34").Value = 7161.357
$ws_CRP.Range("I34").Value = 3592.5557
$ws_CRP.Range("K34").Value = 3592.5557
$ws_CRP.Range("M34").Value = -3390.5557

# CRP row 44 (key=1850)
$ws_CRP.Range("H44").Value = 20000
$ws_CRP.Range("I44").Value = 20000
$ws_CRP.Range("K44").Value = 20000
$ws_CRP.Range("M44").Value = -19558

# CRP row 58 (key=44021)
$ws_CRP.Range("H58").Value = 3889.3157
$ws_CRP.Range("I58").Value = 2376.2942
$ws_CRP.Range("K58").Value = 2376.2942
$ws_CRP.Range("M58").Value = -2173.2942

# CRP row 86 (key=12584)
$ws_CRP.Range("H86").Value = 5416.4443
$ws_CRP.Range("I86").Value = 4559.8
$ws_CRP.Range("J86").Value = 6487.25
$ws_CRP.Range("K86").Value = 4559.8
$ws_CRP.Range("L86").Value = 6487.25
$ws_CRP.Range("M86").Value = -3436.8
$ws_CRP.Range("N86").Value = -8733.25

# CRP row 89 (key=12584)
$ws_CRP.Range("H89").Value = 5416.4443
$ws_CRP.Range("I89").Value = 4559.8
$ws_CRP.Range("J89").Value = 6487.25
$ws_CRP.Range("K89").Value = 22799
$ws_CRP.Range("L89").Value = 32436.25
$ws_CRP.Range("M89").Value = -17183
$ws_CRP.Range("N89").Value = -43668.25

# CRP row 132 (key=44019)
$ws_CRP.Range("H132").Value = 3892.8125
$ws_CRP.Range("I132").Value = 2952.3333
$ws_CRP.Range("K132").Value = 8856.999899999999
$ws_CRP.Range("M132").Value = -6326.999899999999

# CRP row 134 (key=44020)
$ws_CRP.Range("H134").Value = 5542.609
$ws_CRP.Range("I134").Value = 4499.1904
$ws_CRP.Range("K134").Value = 13497.5712
$ws_CRP.Range("M134").Value = -10962.5712

# CRP row 136 (key=44021)
$ws_CRP.Range("H136").Value = 3889.3157
$ws_CRP.Range("I136").Value = 2376.2942
$ws_CRP.Range("K136").Value = 7128.882599999999
$ws_CRP.Range("M136").Value = -4578.882599999999

$ws_CUL = $wb.Worksheets.Item("CUL")
# CUL row 23 (key=4858)
$ws_CUL.Range("H23").Value = 130.85715
$ws_CUL.Range("J23").Value = 130.85715
$ws_CUL.Range("L23").Value = 392.57145
$ws_CUL.Range("N23").Value = -862.5714499999999

# CUL row 80 (key=12890)
$ws_CUL.Range("H80").Value = 4665.8335
$ws_CUL.Range("J80").Value = 4999
$ws_CUL.Range("L80").Value = 14997
$ws_CUL.Range("N80").Value = -16869

# CUL row 83 (key=12890)
$ws_CUL.Range("H83").Value = 4665.8335
$ws_CUL.Range("J83").Value = 4999
$ws_CUL.Range("L83").Value = 44991
$ws_CUL.Range("N83").Value = -54351

# CUL row 122 (key=36078)
$ws_CUL.Range("H122").Value = 1224.75
$ws_CUL.Range("J122").Value = 1449.5
$ws_CUL.Range("L122").Value = 13045.5
$ws_CUL.Range("N122").Value = -17945.5

$ws_GSM = $wb.Worksheets.Item("GSM")
# GSM row 101 (key=18513)
$ws_GSM.Range("H101").Value = 64731.668
$ws_GSM.Range("J101").Value = 64731.668
$ws_GSM.Range("L101").Value = 64731.668
$ws_GSM.Range("N101").Value = -71221.66800000001

# GSM row 113 (key=27710)
$ws_GSM.Range("H113").Value = 2949.3809
$ws_GSM.Range("I113").Value = 2060.7144
$ws_GSM.Range("K113").Value = 2060.7144
$ws_GSM.Range("M113").Value = 109.2856000000002

# GSM row 122 (key=36182)
$ws_GSM.Range("H122").Value = 4297.316
$ws_GSM.Range("I122").Value = 3079.1875
$ws_GSM.Range("J122").Value = 10794
$ws_GSM.Range("K122").Value = 9237.5625
$ws_GSM.Range("L122").Value = 32382
$ws_GSM.Range("M122").Value = -6787.5625
$ws_GSM.Range("N122").Value = -37282

# GSM row 126 (key=36184)
$ws_GSM.Range("H126").Value = 4733.1333
$ws_GSM.Range("I126").Value = 2817.9092
$ws_GSM.Range("J126").Value = 10000
$ws_GSM.Range("K126").Value = 8453.7276
$ws_GSM.Range("L126").Value = 30000
$ws_GSM.Range("M126").Value = -5983.7276
$ws_GSM.Range("N126").Value = -34940

# GSM row 132 (key=44008)
$ws_GSM.Range("H132").Value = 9450.286
$ws_GSM.Range("I132").Value = 8792.654
$ws_GSM.Range("K132").Value = 26377.962
$ws_GSM.Range("M132").Value = -23847.962

$ws_LTW = $wb.Worksheets.Item("LTW")
# LTW row 132 (key=44058)
$ws_LTW.Range("H132").Value = 6662.1665
$ws_LTW.Range("J132").Value = 10143
$ws_LTW.Range("L132").Value = 30429
$ws_LTW.Range("N132").Value = -35489

$ws_WVR = $wb.Worksheets.Item("WVR")
# WVR row 126 (key=36210)
$ws_WVR.Range("H126").Value = 1941.4348
$ws_WVR.Range("I126").Value = 1881.7368
$ws_WVR.Range("J126").Value = 2225
$ws_WVR.Range("K126").Value = 5645.2104
$ws_WVR.Range("L126").Value = 6675
$ws_WVR.Range("M126").Value = -3175.2104
$ws_WVR.Range("N126").Value = -11615

Write-Output "All updates applied."